# 🚌 141: 30/12 19:59 LP1912+6203+6173
# Append newly-scraped rows to each of the 3 schedule sheets, and refresh
# the "Última actualización" / "Total filas" header cells (rows 2 & 3,
# column A) on every sheet.

$wb = $excel.ActiveWorkbook

$stamp = "Última actualización: 30/12/2025 16:59:17"

# ---------------------------------------------------------------------------
# Sheet "LP1912" (columns: A=(blank), B=Hora_Scrap, C=Hora_Llegada,
#                 D=Línea, E=Minutos, F=Parada, G=Fecha)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = $stamp
$ws1.Cells.Item(3, 1).Value = "Total filas: 453"

$data1 = @(
    @(432, "16:59:06", "17:01", "16_SANTA ANA", 2, "LP1912", "30/12/2025"),
    @(433, "16:59:06", "17:05", "11_ETCHEVERRY", 6, "LP1912", "30/12/2025"),
    @(434, "16:59:06", "17:05", "23_HERNANDEZ", 6, "LP1912", "30/12/2025"),
    @(435, "16:59:06", "17:10", "10_OLMOS", 11, "LP1912", "30/12/2025"),
    @(436, "16:59:06", "17:21", "26_HERNANDEZ", 22, "LP1912", "30/12/2025"),
    @(437, "16:59:06", "17:22", "10_OLMOS", 23, "LP1912", "30/12/2025"),
    @(438, "16:59:06", "17:24", "84_COLONIA URQUIZA-ESC 49", 25, "LP1912", "30/12/2025"),
    @(439, "16:59:06", "17:29", "14_ABASTO", 30, "LP1912", "30/12/2025"),
    @(440, "16:59:06", "17:31", "15_ABASTO", 32, "LP1912", "30/12/2025"),
    @(441, "16:59:06", "17:35", "23_HERNANDEZ", 36, "LP1912", "30/12/2025"),
    @(442, "16:59:06", "17:37", "27_EL RETIRO", 38, "LP1912", "30/12/2025"),
    @(443, "16:59:06", "17:38", "17_ROMERO", 39, "LP1912", "30/12/2025"),
    @(444, "16:59:06", "17:41", "16_SANTA ANA", 42, "LP1912", "30/12/2025"),
    @(445, "16:59:06", "17:45", "15_ABASTO", 46, "LP1912", "30/12/2025"),
    @(446, "16:59:06", "17:51", "16_P MOR-167 Y 521", 52, "LP1912", "30/12/2025"),
    @(447, "16:59:06", "17:52", "81_EL PELIGRO", 53, "LP1912", "30/12/2025"),
    @(448, "16:59:06", "18:09", "23_HERNANDEZ", 70, "LP1912", "30/12/2025"),
    @(449, "16:59:06", "18:12", "16_SANTA ANA", 73, "LP1912", "30/12/2025"),
    @(450, "16:59:06", "18:16", "15_ABASTO", 77, "LP1912", "30/12/2025"),
    @(451, "16:59:06", "18:21", "26_HERNANDEZ", 82, "LP1912", "30/12/2025"),
    @(452, "16:59:06", "18:25", "14_ABASTO", 86, "LP1912", "30/12/2025"),
    @(453, "16:59:06", "18:28", "215C_EL PATO", 89, "LP1912", "30/12/2025"),
    @(454, "16:59:06", "18:32", "11X44_ETCHEVERRY", 93, "LP1912", "30/12/2025")
)

foreach ($row in $data1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215" (columns: A=(blank), B=Fecha, C=Hora_Scrap,
#                     D=Hora_Llegada, E=Línea, F=Minutos, G=Parada)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = $stamp
$ws2.Cells.Item(3, 1).Value = "Total filas: 30"

$data2 = @(
    ,@(31, "30/12/2025", "16:59:06", "18:28", "215C_EL PATO", 89, "LP1912")
)

foreach ($row in $data2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173" (columns: A=(blank), B=Fecha, C=Hora_Scrap,
#                    D=Hora_Llegada, E=Línea, F=Minutos, G=Parada)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = $stamp
$ws3.Cells.Item(3, 1).Value = "Total filas: 60"

$data3 = @(
    @(60, "30/12/2025", "16:59:17", "17:29", "215A_LA PLATA", 30, "L6173"),
    @(61, "30/12/2025", "16:59:11", "18:04", "215C_LA PLATA", 65, "L6203")
)

foreach ($row in $data3) {
    $r = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $ws3.Cells.Item($r, 6).Value = $row[5]
    $ws3.Cells.Item($r, 7).Value = $row[6]
}
